# Generate Report for Handoff
# Updates the localization-status report after a handoff generation run:
#  - Overview sheet: refresh "Latest HO Xliff Generate Date" for the
#    handed-off rows.
#  - zh-cn / de-de sheets: refresh "Latest Handoff Datetime" and set
#    "Priority" to "ht" for the same rows.

$wb = $excel.ActiveWorkbook

$rows = @(7, 8, 9, 10, 13, 14)

# --- Overview sheet: column G = "Latest HO Xliff Generate Date" ---
$overview = $wb.Worksheets.Item("Overview")
foreach ($r in $rows) {
    $overview.Range("G$r").Value = "2016-08-22 22:21:43"
}

# --- zh-cn sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$zhcn = $wb.Worksheets.Item("zh-cn")
foreach ($r in $rows) {
    $zhcn.Range("H$r").Value = "2016-08-22 22:21:37"
    $zhcn.Range("E$r").Value = "ht"
}

# --- de-de sheet: column H = "Latest Handoff Datetime", column E = "Priority" ---
$dede = $wb.Worksheets.Item("de-de")
foreach ($r in $rows) {
    $dede.Range("H$r").Value = "2016-08-22 22:21:43"
    $dede.Range("E$r").Value = "ht"
}
